$d = $word.ActiveDocument

$p7xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="6A838310" w14:textId="4024A028" w:rsidR="008D7E50" w:rsidRDefault="0081171E" w:rsidP="008D7E50"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The man needs to transport all three (parrot, cat, seed) to the other side of the river but he can only take one at a time. If</w:t></w:r><w:r w:rsidR="00283A99"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> the man</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> would put the parrot on my shoulder strap the seed to my back and put the cat in the spot that only has room for one more.</w:t></w:r><w:r w:rsidR="00A15948"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> The man’s goal is to get everything to the other side of the river</w:t></w:r><w:r w:rsidR="008D7E50"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidR="008D7E50"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The man only has room for one animal or bag of seed.</w:t></w:r><w:r w:rsidR="00C565CA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00971846"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Sub goals are if the man leaves the parrot with the seed will he eat the seed or if he leaves the cat with the parrot will the cat eat the parrot.</w:t></w:r><w:r w:rsidR="00AE3E4D"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> The man can leave the cat with the seed for a possible solution.</w:t></w:r><w:r w:rsidR="0036140E"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">  There is only one solution I can think of that requires taking one passenger. That is leaving the cat with the seed and taking the parrot across the river first.</w:t></w:r><w:r w:rsidR="003E757F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> The man could try taking two or even three at once if it’s possible.</w:t></w:r><w:r w:rsidR="00283A99"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> If it was </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>me</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> I would try to take at least two or three at once. It would require </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>less</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> trips and save time. I diagnosed a good resolution in my second sentence of this paragraph. If the </w:t></w:r><w:r w:rsidR="00283A99"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve">man strapped the seed to his back, carried the parrot on his shoulder then place the cat in the available space for one passenger then that would work. He could even leave the cat and take the other two. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target = $null
foreach ($p in $d.Paragraphs) { if ($p.Range.Text -like '*If it was me I would*') { $target = $p; break } }
if ($target -eq $null) { throw 'paragraph not found: If it was me I would' }
$target.Range.InsertXML($p7xml)

$p10xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="23BA7834" w14:textId="43A225E0" w:rsidR="000C1374" w:rsidRPr="00DF6009" w:rsidRDefault="00DF6009" w:rsidP="00DF6009"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="00DF6009"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">There are 20 socks in a drawer: 5 pairs of black socks, 3 pairs of brown and 2 pairs of white. You select the socks in the dark and can check them only after a selection has been made. What is the smallest number of socks you need to select to guarantee getting the </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>following:</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target = $null
foreach ($p in $d.Paragraphs) { if ($p.Range.Text -like '*There are 20 socks in a drawer*') { $target = $p; break } }
if ($target -eq $null) { throw 'paragraph not found: There are 20 socks in a drawer' }
$target.Range.InsertXML($p10xml)

$p13xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="3579C3CE" w14:textId="2851C05A" w:rsidR="00DF6009" w:rsidRPr="00DF6009" w:rsidRDefault="000C1374" w:rsidP="00DF6009"><w:pPr><w:pStyle w:val="NormalWeb"/><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>I need to get one matching pair and one matching pair of each color from picking out the smallest number of socks to reach those results.</w:t></w:r><w:r w:rsidR="007059CF"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Its possible to get all three colors in a smaller selection of socks but to guarantee results</w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>;</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> then you would have to gather all 20 socks.</w:t></w:r><w:r w:rsidR="00A04C3B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Overall constraints are</w:t></w:r><w:r w:rsidR="00AF220A"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> picking the socks out in the dark.</w:t></w:r><w:r w:rsidR="00BA3BEA"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Sub goal would be picking out the pairs with a smaller selection.</w:t></w:r><w:r w:rsidR="00A04C3B"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> Picking out eight socks should get you one matching pair but to get a matching pair of each you would have to at least pick 16. </w:t></w:r><w:r w:rsidR="00277CE2"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>The first and second solution meet the goals but still have a chance to fail. It’s not a 100% guarantee.</w:t></w:r><w:r w:rsidR="009F2456"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> There is nothing to say I have to pick out the socks in the dark. If I remove that constraint then I can easily pick out one matching pair and three matching pairs of each color by choosing eight socks. </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target = $null
foreach ($p in $d.Paragraphs) { if ($p.Range.Text -like '*Its possible to get all three colors*') { $target = $p; break } }
if ($target -eq $null) { throw 'paragraph not found: Its possible to get all three colors' }
$target.Range.InsertXML($p13xml)

$p17xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w14:paraId="616F4A46" w14:textId="207CA18E" w:rsidR="00DF6009" w:rsidRPr="00BE242B" w:rsidRDefault="00992447" w:rsidP="00BE242B"><w:pPr><w:spacing w:line="480" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t>The girl is counting on her left hand and if she continues t count like this on what finger will she stop?</w:t></w:r><w:r w:rsidR="00E82423"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t xml:space="preserve"> The little girl would not stop counting if she continued in this manner.</w:t></w:r><w:r w:rsidR="009C261F"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t xml:space="preserve"> Overall goal is to figure out what finger she will stop on.</w:t></w:r><w:r w:rsidR="00831FC6"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t xml:space="preserve"> She is only counting on one hand and she does not have a definitive number she is counting too.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman"/></w:rPr><w:t xml:space="preserve"> Sub goal are counting to 10, counting to 100, and counting to 1000, then figuring out what finger she stops on.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target = $null
foreach ($p in $d.Paragraphs) { if ($p.Range.Text -like '*She is only counting on one hand*') { $target = $p; break } }
if ($target -eq $null) { throw 'paragraph not found: She is only counting on one hand' }
$target.Range.InsertXML($p17xml)
